# Adding Header field in to 1PProfile
#
# The STATUS column (L) on the "Search" sheet was filled with "PASS" for
# every test row (L2:L19). This edit clears those values out (leaving the
# header in L1 untouched), and leaves the sheet scrolled/selected so that
# the now-empty STATUS column is in view and selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "PASS" values out of the STATUS column for every data row,
# without touching the header cell (L1) or any other column.
$ws.Range("L2:L19").ClearContents()

# Scroll the view so column H is the left-most visible column, and leave
# the freshly-cleared STATUS column selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
[void]$ws.Range("L2:L19").Select()
